# Updates the cryptocurrency price/volume table (columns B-E, rows 2-51)
# to match the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.984.54'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.827.04'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").Value = "'" + '1.006'
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").Value = "'" + '311.19'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").Value = "'" + '1.004'
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'" + '0.3696'
$ws.Range("E8").Value = '  +1.56%  '
$ws.Range("D9").Value = "'" + '0.07339'
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = "'" + '0.8753'
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("D11").Value = "'" + '0.07896'
$ws.Range("E11").Value = '  +3.88%  '
$ws.Range("D12").Value = "'" + '19.68'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = '1.839.03'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").Value = "'" + '5.334'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").Value = "'" + '6.554'
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").Value = "'" + '91.40'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = "'" + '0.000008832'
$ws.Range("E18").Value = '  +2.35%  '
$ws.Range("D19").Value = "'" + '1.004'
$ws.Range("E19").Value = '  -0.49%  '
$ws.Range("D20").Value = "'" + '14.78'
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").Value = '27.009.51'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = "'" + '5.098'
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").Value = '2.077.50'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").Value = "'" + '152.09'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").Value = "'" + '1.859'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("D27").Value = "'" + '18.37'
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("D28").Value = "'" + '2.036'
$ws.Range("E28").Value = '  -2.84%  '
$ws.Range("D29").Value = "'" + '5.106'
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").Value = "'" + '115.63'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Value = "'" + '0.08880'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = "'" + '2.959'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = "'" + '0.7302'
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("D34").Value = "'" + '4.437'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").Value = "'" + '1.131'
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("D36").Value = "'" + '2.472'
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").Value = "'" + '1.074'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = "'" + '0.01942'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("D39").Value = "'" + '0.05219'
$ws.Range("E39").Value = '  -1.11%  '
$ws.Range("D40").Value = "'" + '2.960'
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = "'" + '7.088'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("D42").Value = "'" + '0.5151'
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("D43").Value = "'" + '0.1626'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("D44").Value = "'" + '8.150'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").Value = "'" + '0.4823'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'" + '10.20'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = "'" + '1.004'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = "'" + '101.79'
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").Value = "'" + '1.623'
$ws.Range("E49").Value = '  -0.65%  '
$ws.Range("D50").Value = "'" + '0.06207'
$ws.Range("D51").Value = "'" + '64.76'
